# Update "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps as part of regenerating the
# handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G, row 2)
$wsOverview.Range("G2").Value = "2017-03-02 09:01:14"

# zh-cn sheet: "Correspond Handoff Datetime" (column H) and
# "Correspond Handback DateTime" (column L)
$wsZhCn.Range("H2").Value = "2017-03-02 09:00:57"
$wsZhCn.Range("L2").Value = "2017-03-02 09:01:57"

# de-de sheet: "Correspond Handoff Datetime" (column H) and
# "Correspond Handback DateTime" (column L)
$wsDeDe.Range("H2").Value = "2017-03-02 09:01:14"
$wsDeDe.Range("L2").Value = "2017-03-02 09:02:20"
